$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''38.896.46'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.70%  '
$ws.Range('D3').Value = '''2.224.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -6.76%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''296.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.56%  '
$ws.Range('D6').Value = '''80.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -9.27%  '
$ws.Range('D7').Value = '''0.505'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.76%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.457'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.32%  '
$ws.Range('D10').Value = '''0.0769'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.71%  '
$ws.Range('D11').Value = '''27.83'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -10.78%  '
$ws.Range('D12').Value = '''46.00'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -13.40%  '
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('D14').Value = '''2.571.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.63%  '
$ws.Range('D15').Value = '''6.08'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.26%  '
$ws.Range('D16').Value = '''13.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.78%  '
$ws.Range('D17').Value = '''2.239.02'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.02%  '
$ws.Range('D18').Value = '''0.711'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.57%  '
$ws.Range('D19').Value = '''38.863.45'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.54%  '
$ws.Range('D20').Value = '''0.0₃0856'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.11%  '
$ws.Range('E21').Value = '  -7.54%  '
$ws.Range('D22').Value = '''64.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.72%  '
$ws.Range('D23').Value = '''9.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -9.34%  '
$ws.Range('D24').Value = '''224.69'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.77%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  -10.47%  '
$ws.Range('E27').Value = '  -6.11%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '''2.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.61%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''22.07'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.74%  '
$ws.Range('D30').Value = '''8.86'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.83%  '
$ws.Range('D31').Value = '''148.67'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.65%  '
$ws.Range('D32').Value = '''30.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.86%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '''4.77'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.01%  '
$ws.Range('D35').Value = '''2.32'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.70%  '
$ws.Range('D36').Value = '''0.0684'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.58%  '
$ws.Range('D37').Value = '''0.108'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.60%  '
$ws.Range('D38').Value = '''2.66'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.78%  '
$ws.Range('D39').Value = '''0.0949'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.99%  '
$ws.Range('D40').Value = '''14.38'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.25%  '
$ws.Range('D41').Value = '''1.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.42%  '
$ws.Range('D42').Value = '''3.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.92%  '
$ws.Range('D43').Value = '''1.896.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.26%  '
$ws.Range('E44').Value = '  -9.09%  '
$ws.Range('D45').Value = '''0.0252'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.17%  '
$ws.Range('D46').Value = '''16.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.23%  '
$ws.Range('E47').Value = '  -4.22%  '
$ws.Range('D48').Value = '''2.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -10.80%  '
$ws.Range('D49').Value = '''2.440.42'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.67%  '
$ws.Range('D50').Value = '''68.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.32%  '
$ws.Range('D51').Value = '''87.22'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.08%  '
